# Applies the "Groups expenses and split" edit:
#  - Adds a "Checked" column (D) marker to most existing endpoint rows
#  - Adds a new "Financial Goals" endpoint block (rows 41-45)
#  - Adds a new "Groups" endpoint block (rows 47-53)
#  - Adjusts column widths, zoom, and selection to match the authored view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Mark existing rows as "Checked" in column D (style: bold, green font)
# ---------------------------------------------------------------------
$checkedRows = @(3,4,5,6,7,10,11,12,13,14,17,18,19,20,21,24,25,26,27,28,30,31,32,33,34,38)
foreach ($r in $checkedRows) {
    $c = $ws.Cells.Item($r, 4)
    $c.Value = "Checked"
    $c.Font.Bold = $true
    $c.Font.ThemeColor = 10
}

# Style-only placeholder cells in column D (no value)
$ws.Cells.Item(39, 4).Font.Bold = $true
$ws.Cells.Item(40, 4).Font.Bold = $true

Write-Host "step1 done"

# ---------------------------------------------------------------------
# 2) New "Financial Goals" endpoint rows (41-45)
# ---------------------------------------------------------------------
$goalsRows = @(
    @(41, "/api/v1/finance/goals",     "GET",    "List all goals of specific user"),
    @(42, "/api/v1/finance/goals",     "POST",   "Create a new goal"),
    @(43, "/api/v1/finance/goals/{id}","GET",    "Retrieve details of a specific Goal"),
    @(44, "/api/v1/finance/goals/{id}","PUT",    "Edit goal"),
    @(45, "/api/v1/finance/goals/{id}","DELETE", "Delete a goal")
)
foreach ($row in $goalsRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $d = $ws.Cells.Item($r, 4)
    $d.Value = "Checked"
    $d.Font.Bold = $true
    $d.Font.ThemeColor = 10
}

Write-Host "step2 done"

# ---------------------------------------------------------------------
# 3) New "Groups" section header (row 47)
# ---------------------------------------------------------------------
$h = $ws.Cells.Item(47, 1)
$h.Value = "Groups``"
$h.Font.Bold = $true

$h = $ws.Cells.Item(47, 2)
$h.Value = "Methods"
$h.Font.Bold = $true

$h = $ws.Cells.Item(47, 3)
$h.Value = "Description"
$h.Font.Bold = $true

Write-Host "step3 done"

# ---------------------------------------------------------------------
# 4) New "Groups" endpoint rows (48-53), wrap-text body style
# ---------------------------------------------------------------------
$groupsRows = @(
    @(48, "/api/v1/finance/groups/",               "GET",    "List all groups the authenticated user is a member of.",                         $true),
    @(49, "/api/v1/finance/groups/",                "POST",   "Create a new group (user creating is automatically set as admin).",              $true),
    @(50, "/api/v1/finance/groups/{id}/",           "GET",    "Retrieve details of a specific group, including members and expenses.",          $true),
    @(51, "/api/v1/finance/groups/{id}/",           "PUT",    "Update details of a specific group (e.g., name, description).",                  $false),
    @(52, "/api/v1/finance/groups/{id}/",           "DELETE", "Delete a specific group.",                                                        $false),
    @(53, "/api/v1/finance/groups/{id}/add-member/","PATCH",  "Add a user to the group by username (validates if user is already a member).",  $false)
)
foreach ($row in $groupsRows) {
    $r = $row[0]

    $a = $ws.Cells.Item($r, 1)
    $a.Value = $row[1]
    $a.Font.Name = "Arial Unicode MS"
    $a.Font.Size = 10
    $a.WrapText = $true
    $a.VerticalAlignment = -4108

    $b = $ws.Cells.Item($r, 2)
    $b.Value = $row[2]
    $b.WrapText = $true
    $b.VerticalAlignment = -4108

    $c = $ws.Cells.Item($r, 3)
    $c.Value = $row[3]
    $c.WrapText = $true
    $c.VerticalAlignment = -4108

    if ($row[4]) {
        $d = $ws.Cells.Item($r, 4)
        $d.Value = "Checked"
        $d.Font.Bold = $true
        $d.Font.ThemeColor = 10
    }
}

# Rows whose wrapped description needs two visual lines
$ws.Rows.Item(50).RowHeight = 28.8
$ws.Rows.Item(53).RowHeight = 28.8

Write-Host "step4 done"
